# duy lay tien loi ngay 02
# Adds a "6tr interest" entry + a "10tr loan" repayment row at the bottom of
# the "CÔ DIỄM" sheet, and flags an unfinished paperwork note (red highlight)
# next to the interest tracker in row 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Select()

# --- Row 26: interest-tracker entry for 2024-02-03 (serial 45325) plus a
#     red-highlighted note that the paperwork isn't done yet.
$ws.Range("H26").Value = 45325
$ws.Range("I26").Value = 10
$ws.Range("K26").Value = "chưa làm giấy"
$ws.Range("K26").Interior.Color = 255

# --- New ledger rows 76-77: Duy takes 6tr interest on 2023-02-03 (serial
#     44960), then the recurring 10tr loan-to-cô-Diễm repayment line.
$ws.Range("A76").Value = 44960
$ws.Range("B76").Value = "Duy lấy tiền lời 6tr"
$ws.Range("C76").Value = 6000
$ws.Range("D76").Formula = "=D75+C76"

$ws.Range("B77").Value = "Duy cho cô Diễm vay 10tr"
$ws.Range("C77").Value = -10000
$ws.Range("D77").Formula = "=D76+C77"

# --- View state: unfreeze back to the top and leave selection on K35.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("K35").Select()
